$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank log entries for rows 15 and 16
# (Date 2017-11-16 / serial 43055, with Start/End times; Work Time
# is a pre-existing shared formula =ABS(C-B) that recalculates automatically).
$ws.Range("A15").Value = 43055
$ws.Range("B15").Value = 0.65972222222222221
$ws.Range("C15").Value = 0.76388888888888884

$ws.Range("A16").Value = 43055
$ws.Range("B16").Value = 0.78472222222222221
$ws.Range("C16").Value = 0.86805555555555547

# Move the active selection to match the author's final cursor position
$ws.Range("G17").Select()
